# Add output row for employee 7839 (KING), the president with no manager.
# This row is inserted before the existing row 10 (empno 7844, TURNER),
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10, pushing rows 10-14 down to 11-15.
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with KING's data.
# empno
$ws.Range("A10").Value = 7839
# emp_name
$ws.Range("B10").Value = "KING"
# mgr_name intentionally left blank - KING is the top of the hierarchy and
# has no manager.
